$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4375153333333333
$ws.Range("H2").Value = 1.312546
$ws.Range("I2").Value = 0.002535486401940996
$ws.Range("J2").Value = 0.002555908833496712
$ws.Range("M2").Value = 1.428094
$ws.Range("N2").Value = 4.284282
$ws.Range("O2").Value = 0.01363994948788035
$ws.Range("P2").Value = 0.01373582358494966
$ws.Range("Q2").Value = 0.6248130224413333
$ws.Range("R2").Value = 5.623317201972
$ws.Range("S2").Value = 0.00003458390644968269
$ws.Range("T2").Value = 0.00003510751283612532

$ws.Range("G3").Value = 0.4375153333333333
$ws.Range("H3").Value = 1.312546
$ws.Range("I3").Value = 0.002535486401940996
$ws.Range("J3").Value = 0.002555908833496712
$ws.Range("M3").Value = 2.242972
$ws.Range("N3").Value = 6.728916
$ws.Range("O3").Value = 0.02142297690679323
$ws.Range("P3").Value = 0.02157355727142731
$ws.Range("Q3").Value = 0.9813346422373332
$ws.Range("R3").Value = 8.832011780136
$ws.Range("S3").Value = 0.00005431766663627022
$ws.Range("T3").Value = 0.00005514004559998829

$ws.Range("G4").Value = 0.4375153333333333
$ws.Range("H4").Value = 1.312546
$ws.Range("I4").Value = 0.002535486401940996
$ws.Range("J4").Value = 0.002555908833496712
$ws.Range("M4").Value = 56.98919799999999
$ws.Range("N4").Value = 170.967594
$ws.Range("O4").Value = 0.5443127567756828
$ws.Range("P4").Value = 0.5481386869322091
$ws.Range("Q4").Value = 24.93364795936933
$ws.Range("R4").Value = 224.402831634324
$ws.Range("S4").Value = 0.001380097593207761
$ws.Range("T4").Value = 0.001400992511911322

$ws.Range("G5").Value = 0.4375153333333333
$ws.Range("H5").Value = 1.312546
$ws.Range("I5").Value = 0.002535486401940996
$ws.Range("J5").Value = 0.002555908833496712
$ws.Range("M5").Value = 2.19236
$ws.Range("N5").Value = 4.38472
$ws.Range("O5").Value = 0.0209395737670275
$ws.Range("P5").Value = 0.01405783755350383
$ws.Range("Q5").Value = 0.9591911161866665
$ws.Range("R5").Value = 5.75514669712
$ws.Range("S5").Value = 0.00005309200454873862
$ws.Range("T5").Value = 0.00003593055118286223

$ws.Range("G6").Value = 0.4375153333333333
$ws.Range("H6").Value = 1.312546
$ws.Range("I6").Value = 0.002535486401940996
$ws.Range("J6").Value = 0.002555908833496712
$ws.Range("M6").Value = 41.84673733333333
$ws.Range("N6").Value = 125.540212
$ws.Range("O6").Value = 0.399684743062616
$ws.Range("P6").Value = 0.4024940946579102
$ws.Range("Q6").Value = 18.30858923330577
$ws.Range("R6").Value = 164.777303099752
$ws.Range("S6").Value = 0.001013395231098544
$ws.Range("T6").Value = 0.001028738211966414

$ws.Range("G7").Value = 1.004357666666667
$ws.Range("H7").Value = 3.013073
$ws.Range("I7").Value = 0.005820447907772805
$ws.Range("J7").Value = 0.005867329523437988
$ws.Range("M7").Value = 1.428094
$ws.Range("N7").Value = 4.284282
$ws.Range("O7").Value = 0.01363994948788035
$ws.Range("P7").Value = 0.01373582358494966
$ws.Range("Q7").Value = 1.434317157620667
$ws.Range("R7").Value = 12.908854418586
$ws.Range("S7").Value = 0.00007939061545885994
$ws.Range("T7").Value = 0.000080592603248711

$ws.Range("G8").Value = 1.004357666666667
$ws.Range("H8").Value = 3.013073
$ws.Range("I8").Value = 0.005820447907772805
$ws.Range("J8").Value = 0.005867329523437988
$ws.Range("M8").Value = 2.242972
$ws.Range("N8").Value = 6.728916
$ws.Range("O8").Value = 0.02142297690679323
$ws.Range("P8").Value = 0.02157355727142731
$ws.Range("Q8").Value = 2.252746124318667
$ws.Range("R8").Value = 20.274715118868
$ws.Range("S8").Value = 0.0001246913211154098
$ws.Range("T8").Value = 0.0001265791695042258

$ws.Range("G9").Value = 1.004357666666667
$ws.Range("H9").Value = 3.013073
$ws.Range("I9").Value = 0.005820447907772805
$ws.Range("J9").Value = 0.005867329523437988
$ws.Range("M9").Value = 56.98919799999999
$ws.Range("N9").Value = 170.967594
$ws.Range("O9").Value = 0.5443127567756828
$ws.Range("P9").Value = 0.5481386869322091
$ws.Range("Q9").Value = 57.23753792848466
$ws.Range("R9").Value = 515.137841356362
$ws.Range("S9").Value = 0.003168144046349071
$ws.Range("T9").Value = 0.003216110300775883

$ws.Range("G10").Value = 1.004357666666667
$ws.Range("H10").Value = 3.013073
$ws.Range("I10").Value = 0.005820447907772805
$ws.Range("J10").Value = 0.005867329523437988
$ws.Range("M10").Value = 2.19236
$ws.Range("N10").Value = 4.38472
$ws.Range("O10").Value = 0.0209395737670275
$ws.Range("P10").Value = 0.01405783755350383
$ws.Range("Q10").Value = 2.201913574093333
$ws.Range("R10").Value = 13.21148144456
$ws.Range("S10").Value = 0.0001218776983219495
$ws.Range("T10").Value = 0.00008248196531336827

$ws.Range("G11").Value = 1.004357666666667
$ws.Range("H11").Value = 3.013073
$ws.Range("I11").Value = 0.005820447907772805
$ws.Range("J11").Value = 0.005867329523437988
$ws.Range("M11").Value = 41.84673733333333
$ws.Range("N11").Value = 125.540212
$ws.Range("O11").Value = 0.399684743062616
$ws.Range("P11").Value = 0.4024940946579102
$ws.Range("Q11").Value = 42.02909146571955
$ws.Range("R11").Value = 378.261823191476
$ws.Range("S11").Value = 0.002326344226527515
$ws.Range("T11").Value = 0.002361565484595801

$ws.Range("G12").Value = 99.58055866666666
$ws.Range("H12").Value = 298.741676
$ws.Range("I12").Value = 0.577088694179909
$ws.Range("J12").Value = 0.5817369361698658
$ws.Range("M12").Value = 1.428094
$ws.Range("N12").Value = 4.284282
$ws.Range("O12").Value = 0.01363994948788035
$ws.Range("P12").Value = 0.01373582358494966
$ws.Range("Q12").Value = 142.2103983485146
$ws.Range("R12").Value = 1279.893585136632
$ws.Range("S12").Value = 0.007871460638640792
$ws.Range("T12").Value = 0.0079906359280784

$ws.Range("G13").Value = 99.58055866666666
$ws.Range("H13").Value = 298.741676
$ws.Range("I13").Value = 0.577088694179909
$ws.Range("J13").Value = 0.5817369361698658
$ws.Range("M13").Value = 2.242972
$ws.Range("N13").Value = 6.728916
$ws.Range("O13").Value = 0.02142297690679323
$ws.Range("P13").Value = 0.02157355727142731
$ws.Range("Q13").Value = 223.3564048336906
$ws.Range("R13").Value = 2010.207643503216
$ws.Range("S13").Value = 0.01236295776858765
$ws.Range("T13").Value = 0.01255013510936526

$ws.Range("G14").Value = 99.58055866666666
$ws.Range("H14").Value = 298.741676
$ws.Range("I14").Value = 0.577088694179909
$ws.Range("J14").Value = 0.5817369361698658
$ws.Range("M14").Value = 56.98919799999999
$ws.Range("N14").Value = 170.967594
$ws.Range("O14").Value = 0.5443127567756828
$ws.Range("P14").Value = 0.5481386869322091
$ws.Range("Q14").Value = 5675.016174805281
$ws.Range("R14").Value = 51075.14557324754
$ws.Range("S14").Value = 0.3141167380331452
$ws.Range("T14").Value = 0.3188725203321166

$ws.Range("G15").Value = 99.58055866666666
$ws.Range("H15").Value = 298.741676
$ws.Range("I15").Value = 0.577088694179909
$ws.Range("J15").Value = 0.5817369361698658
$ws.Range("M15").Value = 2.19236
$ws.Range("N15").Value = 4.38472
$ws.Range("O15").Value = 0.0209395737670275
$ws.Range("P15").Value = 0.01405783755350383
$ws.Range("Q15").Value = 218.3164335984533
$ws.Range("R15").Value = 1309.89860159072
$ws.Range("S15").Value = 0.01208399128189778
$ws.Range("T15").Value = 0.008177963347548998

$ws.Range("G16").Value = 99.58055866666666
$ws.Range("H16").Value = 298.741676
$ws.Range("I16").Value = 0.577088694179909
$ws.Range("J16").Value = 0.5817369361698658
$ws.Range("M16").Value = 41.84673733333333
$ws.Range("N16").Value = 125.540212
$ws.Range("O16").Value = 0.399684743062616
$ws.Range("P16").Value = 0.4024940946579102
$ws.Range("Q16").Value = 4167.12148203059
$ws.Range("R16").Value = 37504.09333827531
$ws.Range("S16").Value = 0.2306535464576376
$ws.Range("T16").Value = 0.2341456814527566

$ws.Range("G17").Value = 4.1363315
$ws.Range("H17").Value = 8.272663
$ws.Range("I17").Value = 0.02397084507248554
$ws.Range("J17").Value = 0.01610928107528529
$ws.Range("M17").Value = 1.428094
$ws.Range("N17").Value = 4.284282
$ws.Range("O17").Value = 0.01363994948788035
$ws.Range("P17").Value = 0.01373582358494966
$ws.Range("Q17").Value = 5.907070197161
$ws.Range("R17").Value = 35.442421182966
$ws.Range("S17").Value = 0.0003269611159705084
$ws.Range("T17").Value = 0.000221274242930487

$ws.Range("G18").Value = 4.1363315
$ws.Range("H18").Value = 8.272663
$ws.Range("I18").Value = 0.02397084507248554
$ws.Range("J18").Value = 0.01610928107528529
$ws.Range("M18").Value = 2.242972
$ws.Range("N18").Value = 6.728916
$ws.Range("O18").Value = 0.02142297690679323
$ws.Range("P18").Value = 0.02157355727142731
$ws.Range("Q18").Value = 9.277675737217999
$ws.Range("R18").Value = 55.666054423308
$ws.Range("S18").Value = 0.0005135268604241761
$ws.Range("T18").Value = 0.0003475344978791875

$ws.Range("G19").Value = 4.1363315
$ws.Range("H19").Value = 8.272663
$ws.Range("I19").Value = 0.02397084507248554
$ws.Range("J19").Value = 0.01610928107528529
$ws.Range("M19").Value = 56.98919799999999
$ws.Range("N19").Value = 170.967594
$ws.Range("O19").Value = 0.5443127567756828
$ws.Range("P19").Value = 0.5481386869322091
$ws.Range("Q19").Value = 235.726214847137
$ws.Range("R19").Value = 1414.357289082822
$ws.Range("S19").Value = 0.0130476367636474
$ws.Range("T19").Value = 0.008830120176028765

$ws.Range("G20").Value = 4.1363315
$ws.Range("H20").Value = 8.272663
$ws.Range("I20").Value = 0.02397084507248554
$ws.Range("J20").Value = 0.01610928107528529
$ws.Range("M20").Value = 2.19236
$ws.Range("N20").Value = 4.38472
$ws.Range("O20").Value = 0.0209395737670275
$ws.Range("P20").Value = 0.01405783755350383
$ws.Range("Q20").Value = 9.06832772734
$ws.Range("R20").Value = 36.27331090936
$ws.Range("S20").Value = 0.0005019392786532986
$ws.Range("T20").Value = 0.0002264616564600941

$ws.Range("G21").Value = 4.1363315
$ws.Range("H21").Value = 8.272663
$ws.Range("I21").Value = 0.02397084507248554
$ws.Range("J21").Value = 0.01610928107528529
$ws.Range("M21").Value = 41.84673733333333
$ws.Range("N21").Value = 125.540212
$ws.Range("O21").Value = 0.399684743062616
$ws.Range("P21").Value = 0.4024940946579102
$ws.Range("Q21").Value = 173.0919778040926
$ws.Range("R21").Value = 1038.551866824556
$ws.Range("S21").Value = 0.00958078105379016
$ws.Range("T21").Value = 0.006483890501986759

$ws.Range("G22").Value = 67.39800266666667
$ws.Range("H22").Value = 202.194008
$ws.Range("I22").Value = 0.3905845264378918
$ws.Range("J22").Value = 0.3937305443979143
$ws.Range("M22").Value = 1.428094
$ws.Range("N22").Value = 4.284282
$ws.Range("O22").Value = 0.01363994948788035
$ws.Range("P22").Value = 0.01373582358494966
$ws.Range("Q22").Value = 96.25068322025066
$ws.Range("R22").Value = 866.256148982256
$ws.Range("S22").Value = 0.005327553211360512
$ws.Range("T22").Value = 0.005408213297855942

$ws.Range("G23").Value = 67.39800266666667
$ws.Range("H23").Value = 202.194008
$ws.Range("I23").Value = 0.3905845264378918
$ws.Range("J23").Value = 0.3937305443979143
$ws.Range("M23").Value = 2.242972
$ws.Range("N23").Value = 6.728916
$ws.Range("O23").Value = 0.02142297690679323
$ws.Range("P23").Value = 0.02157355727142731
$ws.Range("Q23").Value = 151.1718328372587
$ws.Range("R23").Value = 1360.546495535328
$ws.Range("S23").Value = 0.008367483290029726
$ws.Range("T23").Value = 0.008494168449078658

$ws.Range("G24").Value = 67.39800266666667
$ws.Range("H24").Value = 202.194008
$ws.Range("I24").Value = 0.3905845264378918
$ws.Range("J24").Value = 0.3937305443979143
$ws.Range("M24").Value = 56.98919799999999
$ws.Range("N24").Value = 170.967594
$ws.Range("O24").Value = 0.5443127567756828
$ws.Range("P24").Value = 0.5481386869322091
$ws.Range("Q24").Value = 3840.958118775195
$ws.Range("R24").Value = 34568.62306897675
$ws.Range("S24").Value = 0.2126001403393334
$ws.Range("T24").Value = 0.2158189436113766

$ws.Range("G25").Value = 67.39800266666667
$ws.Range("H25").Value = 202.194008
$ws.Range("I25").Value = 0.3905845264378918
$ws.Range("J25").Value = 0.3937305443979143
$ws.Range("M25").Value = 2.19236
$ws.Range("N25").Value = 4.38472
$ws.Range("O25").Value = 0.0209395737670275
$ws.Range("P25").Value = 0.01405783755350383
$ws.Range("Q25").Value = 147.7606851262933
$ws.Range("R25").Value = 886.56411075776
$ws.Range("S25").Value = 0.008178673503605738
$ws.Range("T25").Value = 0.005535000032998506

$ws.Range("G26").Value = 67.39800266666667
$ws.Range("H26").Value = 202.194008
$ws.Range("I26").Value = 0.3905845264378918
$ws.Range("J26").Value = 0.3937305443979143
$ws.Range("M26").Value = 41.84673733333333
$ws.Range("N26").Value = 125.540212
$ws.Range("O26").Value = 0.399684743062616
$ws.Range("P26").Value = 0.4024940946579102
$ws.Range("Q26").Value = 2820.3865143833
$ws.Range("R26").Value = 25383.4786294497
$ws.Range("S26").Value = 0.1561106760935623
$ws.Range("T26").Value = 0.1584742190066046
